$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Issue date" and "Python version" info cells
$ws.Range("A5").Value = "Issue date: 10/12/2020 11:27:59"
$ws.Range("A6").Value = "Python version: Python 3.7.6"

# Update the results table style
$table = $ws.ListObjects.Item(1)
$table.TableStyle = "TableStyleLight10"
